$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '70.773.05'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '3.800.57'
$ws.Range('E3').Value = '  -1.36%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.08%  '
Set-TextValue 'D5' '704.09'
$ws.Range('E5').Value = '  +1.15%  '
Set-TextValue 'D6' '170.36'
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('D7').Value = '3.799.38'
$ws.Range('E7').Value = '  -1.38%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  -1.93%  '
Set-TextValue 'D11' '7.38'
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').Value = '4.437.16'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '3.820.73'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '70.689.20'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('E18').Value = '  -0.01%  '
Set-TextValue 'D19' '7.14'
$ws.Range('E19').Value = '  -1.76%  '
Set-TextValue 'D20' '17.37'
$ws.Range('E20').Value = '  -2.14%  '
Set-TextValue 'D21' '497.08'
$ws.Range('E21').Value = '  +0.18%  '
Set-TextValue 'D22' '10.61'
$ws.Range('E22').Value = '  -5.00%  '
Set-TextValue 'D23' '0.727'
$ws.Range('E23').Value = '  +0.46%  '
Set-TextValue 'D24' '84.83'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('E26').Value = '  -2.00%  '
Set-TextValue 'D27' '10.43'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = '3.947.78'
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  -4.71%  '
Set-TextValue 'D31' '3.08'
$ws.Range('E31').Value = '  -2.61%  '
Set-TextValue 'D32' '7.31'
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('E33').Value = '  -4.15%  '
Set-TextValue 'D34' '29.08'
$ws.Range('E34').Value = '  -2.31%  '
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '3.767.14'
$ws.Range('E37').Value = '  -1.00%  '
Set-TextValue 'D38' '9.08'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('E39').Value = '  -3.79%  '
Set-TextValue 'D41' '2.32'
$ws.Range('E41').Value = '  -2.74%  '
Set-TextValue 'D42' '5.94'
$ws.Range('E42').Value = '  -2.06%  '
$ws.Range('E43').Value = '  -3.98%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('E45').Value = '  +0.02%  '
Set-TextValue 'D46' '0.000321'
$ws.Range('E46').Value = '  +4.63%  '
Set-TextValue 'D47' '164.65'
$ws.Range('E47').Value = '  +0.01%  '
Set-TextValue 'D48' '425.66'
$ws.Range('E48').Value = '  +1.36%  '
Set-TextValue 'D49' '48.75'
$ws.Range('E49').Value = '  +0.15%  '
Set-TextValue 'D50' '8.60'
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('E51').Value = '  -1.75%  '
